$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 special case: coin name and link change (Aave -> RenderToken)
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# Update Price (D) and Volume(1h) (E) columns for each changed row.
# Price values are set via a text number-format to prevent Excel from
# auto-converting strings such as "1.012" or "26.933.35" into numbers,
# then the style is restored to Normal so no extra formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.933.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.786.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.90%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4235"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.20%  "

$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07146"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8409"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.795.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06823"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008680"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.195.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.052"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.078.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.945"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("E27").Value = "  -3.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.009"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.619"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08936"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7205"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.54%  "

$ws.Range("E34").Value = "  -4.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.011"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.084"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.081"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01897"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05077"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4939"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1615"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.506"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.952"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.011"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06278"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4477"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.575"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.695"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.95%  "

